$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Values are written with a leading apostrophe so the engine keeps them as
# literal text (matching the source data which stores numbers/percentages
# as plain strings), then ClearFormats() strips the quote-prefix style that
# the apostrophe trick would otherwise leave behind, restoring the default
# (unstyled) cell formatting.
$ws.Range("D2").Value = "'273.02"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'0.97%"
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'26.77"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'0.24%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'4.906"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'4.15%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'0.06316"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'3.41%"
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'6.909"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'2.52%"
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'3.351"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'5.54%"
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'1.385"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'54.33%"
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.8833"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'3.17%"
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.1475"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'2.95%"
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.05099"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'3.05%"
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.07387"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'4.00%"
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.03177"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'0.00%"
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'0.09044"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'0.20%"
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.001564"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'1.42%"
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'0.0006313"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'3.76%"
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'0.006020"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'0.70%"
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'3.472"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'0.27%"
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'2.283"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'0.86%"
$ws.Range("E19").ClearFormats()
$ws.Range("D21").Value = "'0.1334"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'4.22%"
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'3.927"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'2.19%"
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'0.04336"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'2.54%"
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'0.001176"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'-0.19%"
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'0.003641"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'-12.32%"
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'0.09%"
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'0.0001699"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'1.12%"
$ws.Range("E27").ClearFormats()
$ws.Range("D40").Value = "'0.04055"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'2.77%"
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.006607"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'57.74%"
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.1163"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'3.90%"
$ws.Range("E42").ClearFormats()
$ws.Range("E43").Value = "'9.14%"
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.01258"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'4.72%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.00005335"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'3.89%"
$ws.Range("E45").ClearFormats()
$ws.Range("E46").Value = "'142.27%"
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'0.02122"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'-13.30%"
$ws.Range("E47").ClearFormats()
